$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "ATUPropertiesPath" row -> "Screenshots path" row
$ws.Range("B15").Value = "Screenshots path"
$ws.Range("C15").Value = "D:\\Tookitaki\\test-output"

# Update the active selection shown in the sheet view
$ws.Range("B15").Select()
